$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.028.29"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.635.96"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "'214.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").Value = "'18.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.704.71"
$ws.Range("E12").Value = "  +4.11%  "
$ws.Range("D13").Value = "1.865.49"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "'62.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "26.043.00"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "'190.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "'4.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'143.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("B27").Value = "BinanceUSD"
$ws.Range("C27").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D27").Value = "'1.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "'0.0486"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("D37").Value = "1.130.82"
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("D40").Value = "'0.0156"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "'98.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'0.794"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "'55.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'0.0928"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.26%  "
